# Commit: "Finally got the globe to display! All JSON datasets were required
# to be the same length. Unfortunately this also causes some wonky
# positioning for years that are not 2016."
#
# 1) Data fix: the Krakow entry used for the chart label loses its
#    "(Cracow)" qualifier so the (now padded/aligned) dataset lines up with
#    the other years' JSON.
# 2) View-state fallout: the selection/scroll position saved in the sheet
#    drifted to the bottom of the data (row ~209) instead of the top.
# 3) Minor column-width nudge that came along with the same save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the shared-string label used in column D for the Krakow row.
$ws.Range("D169").Value = "Krakow"

# 2. Update the view state: scroll so row 188 is at the top and the
#    active/selected cell becomes D209.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 188
$ws.Range("D209").Select()

# 3. Slightly narrower default column width (14.4438775510204 -> 14.1734693877551
#    character-width units in the saved OOXML). Apply it across the same
#    column range (1 .. 1025, i.e. A:AMK) the workbook already carries
#    custom width info for.
$ws.Range("A1:AMK1").EntireColumn.ColumnWidth = 13.33
